$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.073.44'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '1.872.35'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '''241.89'
$ws.Range("E5").Value = '  -2.14%  '
$ws.Range("D6").Value = '''0.9997'
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").Value = '''0.4860'
$ws.Range("D8").Value = '''0.2880'
$ws.Range("E8").Value = '  -2.16%  '
$ws.Range("D9").Value = '''0.06553'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").Value = '1.873.51'
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("D11").Value = '''16.27'
$ws.Range("E11").Value = '  -4.38%  '
$ws.Range("D12").Value = '''0.07178'
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("E13").Value = '  -2.54%  '
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '''85.68'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '''4.892'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = '30.022.63'
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("D17").Value = '''1.000'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '''0.000007728'
$ws.Range("E18").Value = '  -3.53%  '
$ws.Range("E19").Value = '  -1.53%  '
$ws.Range("D20").Value = '2.113.75'
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").Value = '''4.725'
$ws.Range("E22").Value = '  -1.24%  '
$ws.Range("D23").Value = '''5.814'
$ws.Range("E23").Value = '  +2.83%  '
$ws.Range("D24").Value = '''9.123'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("D25").Value = '''151.53'
$ws.Range("E25").Value = '  +2.51%  '
$ws.Range("D26").Value = '''142.08'
$ws.Range("E26").Value = '  +5.51%  '
$ws.Range("D27").Value = '''16.87'
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").Value = '''1.865'
$ws.Range("E28").Value = '  -4.15%  '
$ws.Range("D29").Value = '''1.389'
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("D30").Value = '''4.174'
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("D31").Value = '''0.08749'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '''3.963'
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").Value = '''0.05088'
$ws.Range("E33").Value = '  -1.41%  '
$ws.Range("D34").Value = '''0.7077'
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("E35").Value = '  -2.17%  '
$ws.Range("D36").Value = '''2.669'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '''0.01828'
$ws.Range("E37").Value = '  +9.85%  '
$ws.Range("D38").Value = '''2.673'
$ws.Range("E38").Value = '  -3.88%  '
$ws.Range("D39").Value = '''2.131'
$ws.Range("E39").Value = '  -5.19%  '
$ws.Range("D40").Value = '''0.9219'
$ws.Range("E40").Value = '  -2.17%  '
$ws.Range("D41").Value = '''0.9989'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").Value = '''103.54'
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("D43").Value = '''5.730'
$ws.Range("E43").Value = '  -6.27%  '
$ws.Range("D44").Value = '''0.4198'
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").Value = '''7.367'
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("D46").Value = '''0.1272'
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("D47").Value = '''0.05710'
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").Value = '''32.56'
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("D49").Value = '''8.243'
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").Value = '''0.3728'
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.329'
$ws.Range("E51").Value = '  -1.21%  '
